# Fix final affichage mapping posologie
$wb = $excel.ActiveWorkbook

# 1. Update the metadata "Date" value on the Metadata sheet (B8).
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-12-27T10:41:44+00:00"

# 2. Clear the stray "Elément_posologie/Fréquence_structurée/Frq_échelle"
#    values left over in rows 4-9, column A, of "Mapping Table 2" (row 3
#    keeps its value; only the duplicated rows below it are cleared).
$ws = $wb.Worksheets.Item("Mapping Table 2")
$ws.Range("A4:A9").ClearContents()
